# "added 4wk low sales check"
# Updates the forecast re-run results on the "Forecast Comparison" sheet
# (MyForecast, Inventory Coverage, Stockout Risk, Reorder Urgency and
# Seasonality Index columns) and rolls the new totals up into the
# "Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# Row 2 (W10)
$ws.Range("D2").Value = 46
$ws.Range("H2").Value = 8.630000000000001
$ws.Range("L2").Value = 0.92

# Row 3 (W11)
$ws.Range("D3").Value = 46
$ws.Range("H3").Value = 7.63
$ws.Range("L3").Value = 0.99

# Row 4 (W12)
$ws.Range("D4").Value = 47
$ws.Range("H4").Value = 6.49
$ws.Range("L4").Value = 1.11

# Row 5 (W13)
$ws.Range("D5").Value = 47
$ws.Range("H5").Value = 5.49
$ws.Range("L5").Value = 0.83

# Row 6 (W14)
$ws.Range("D6").Value = 47
$ws.Range("H6").Value = 4.49
$ws.Range("L6").Value = 1.05

# Row 7 (W15)
$ws.Range("D7").Value = 48
$ws.Range("H7").Value = 3.42
$ws.Range("L7").Value = 1.11

# Row 8 (W16)
$ws.Range("D8").Value = 48
$ws.Range("H8").Value = 2.42
$ws.Range("L8").Value = 1

# Row 9 (W17)
$ws.Range("D9").Value = 48
$ws.Range("H9").Value = 1.42
$ws.Range("L9").Value = 0.93

# Row 10 (W18)
$ws.Range("D10").Value = 48
$ws.Range("H10").Value = 0.42
$ws.Range("I10").Value = "High"
$ws.Range("J10").Value = "Urgent"
$ws.Range("L10").Value = 0.97

# Row 11 (W19)
$ws.Range("D11").Value = 49
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = "High"
$ws.Range("L11").Value = 1.07

# Row 12 (W20)
$ws.Range("D12").Value = 49
$ws.Range("L12").Value = 0.97

# Row 13 (W21)
$ws.Range("D13").Value = 49
$ws.Range("L13").Value = 1

# Row 14 (W22)
$ws.Range("D14").Value = 49
$ws.Range("L14").Value = 1.03

# Row 15 (W23)
$ws.Range("D15").Value = 49
$ws.Range("L15").Value = 1.03

# Row 16 (W24)
$ws.Range("D16").Value = 50
$ws.Range("L16").Value = 0.8100000000000001

# Row 17 (W25)
$ws.Range("D17").Value = 50
$ws.Range("L17").Value = 1.16

# Roll the refreshed MyForecast column up into the Summary sheet.
# (Leading apostrophe keeps these as text cells, matching the column's
# existing inline-string formatting instead of flipping them to numbers.)
$summary.Range("B9").Value = "'770"
$summary.Range("B10").Value = "'377"
$summary.Range("B11").Value = "'186"
$summary.Range("B12").Value = "'50"
$summary.Range("B14").Value = "'46"
